$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values of the columns that move (D, M, N, O, P, Q, S) for rows 3-9.
# Use .Value2 because .Value does not reliably read through this COM shim.
$snapD = @{}
$snapM = @{}
$snapN = @{}
$snapO = @{}
$snapP = @{}
$snapQ = @{}
$snapS = @{}

foreach ($r in 3..9) {
    $snapD[$r] = $ws.Range("D$r").Value2
    $snapM[$r] = $ws.Range("M$r").Value2
    $snapN[$r] = $ws.Range("N$r").Value2
    $snapO[$r] = $ws.Range("O$r").Value2
    $snapP[$r] = $ws.Range("P$r").Value2
    $snapQ[$r] = $ws.Range("Q$r").Value2
    $snapS[$r] = $ws.Range("S$r").Value2
}

# Mapping of target row -> source row: the values that should end up in the
# target row are the ones currently sitting in the source row.
$mapping = @{
    3 = 9
    4 = 3
    5 = 4
    6 = 8
    7 = 6
    8 = 5
    9 = 7
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]

    $ws.Range("D$target").Value2 = $snapD[$source]
    $ws.Range("M$target").Value2 = $snapM[$source]
    $ws.Range("N$target").Value2 = $snapN[$source]
    $ws.Range("O$target").Value2 = $snapO[$source]
    $ws.Range("P$target").Value2 = $snapP[$source]
    $ws.Range("Q$target").Value2 = $snapQ[$source]
    $ws.Range("S$target").Value2 = $snapS[$source]
}
